# 0.1.2: Client mode is implemented for micronaut.
#
# The resource-bundle sheet ("ja") documents the keys used by the
# blancoRestGeneratorKt code generator. This change adds three new
# documentation rows describing the newly-implemented REST *client*
# executor (interface / request-arg / response-return), inserted right
# after the existing server-side "EXECUTOR" rows (row 89) and before the
# "AUTHFLAG" rows (old row 90, which shifts down to row 93).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ja")

# Insert 3 blank rows before the old row 90 - everything below (old rows
# 90:106) shifts down to 93:109, and Excel re-maps the shared formulas in
# column A plus the dataValidation sqref automatically.
$ws.Rows("90:92").Insert()

# The newly inserted rows don't inherit the numbering-column formatting /
# formula from the row above automatically - copy the look from row 89
# (the last row of the block above) and then write the "+1" running
# counter formula explicitly, matching the existing A83:A106(->A109)
# shared-formula family.
$ws.Range("A89").Copy()
$ws.Range("A90:A92").PasteSpecial(-4122)
$ws.Range("A90").Formula = "=A89+1"
$ws.Range("A91").Formula = "=A90+1"
$ws.Range("A92").Formula = "=A91+1"

# New row 90: XML2SOURCE_FILE.CLIENT.EXECUTOR.DESCRIPTION
$ws.Range("B90").Value = "XML2SOURCE_FILE.CLIENT.EXECUTOR.DESCRIPTION"
$ws.Range("C90").Value = "クライアントとしてAPIよぶためのインタフェイスです"

# New row 91: XML2SOURCE_FILE.CLIENT.EXECUTOR.ARG.LANGDOC
$ws.Range("B91").Value = "XML2SOURCE_FILE.CLIENT.EXECUTOR.ARG.LANGDOC"
$ws.Range("C91").Value = "クライアントとして渡すリクエスト情報です"

# New row 92: XML2SOURCE_FILE.CLIENT.EXECUTOR.RETURN.LANGDOC
$ws.Range("B92").Value = "XML2SOURCE_FILE.CLIENT.EXECUTOR.RETURN.LANGDOC"
$ws.Range("C92").Value = "サーバから戻されるレスポンス情報です"

# Leave the cursor/selection where the author ended up after the edit.
[void]$ws.Range("B112").Select()
